$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.genomeweb.com/cancer/angle-looks-dual-ctc-ctdna-analysis-illumina-tie-expand-cancer-dx-cdx-business"
$keywords = "CDx, ctDNA"
$title = "Angle Looks to Dual CTC, ctDNA Analysis, Illumina Tie-up to Expand Cancer Dx, CDx Business"

$row = 14

# Add the new link (column A) as a hyperlink, matching the style used by the
# other link cells in the column.
$ws.Hyperlinks.Add($ws.Range("A$row"), $url)
$ws.Range("A$row").Style = $ws.Range("A13").Style

# Fill in the keywords/title columns for the new row.
$ws.Range("B$row").Value = $keywords
$ws.Range("C$row").Value = $title
